$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / label text updates (shared strings) ---
$ws.Range("C1").Value = "Integral"
$ws.Range("D1").Value = "Time"
$ws.Range("A13").Value = "Avg"

# --- New dataset (rows 2-11) for columns B (time step), C (integral-ish), D (time) ---
$bVals = @(0.001004389332,0.001042091696,0.0009909262079,0.00098704745549999993,0.0010614461400000001,0.00099417321120000009,0.0010399842750000001,0.001329988666,0.0010859592130000001,0.00097087162040000003)
$cVals = @(0.19343095888967901,0.19493938548356601,0.193268569682986,0.19249050981672899,0.19398757858678101,0.19150072072296401,0.19405432266606101,0.19398518120632,0.19455134277262301,0.192046511471413)
$dVals = @(0.184820287,0.16273359500000001,0.156800369,0.14464328300000001,0.13743134200000001,0.120320991,0.12225472,0.111936202,0.119138261,0.114771111)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
}

# --- Summary rows: re-enter as independent (non-shared) formulas ---
$ws.Range("B13").Formula = "=AVERAGE(B2:B11)"
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"
$ws.Range("D13").Formula = "=AVERAGE(D2:D11)"

$ws.Range("B14").Formula = "=_xlfn.STDEV.S(B2:B11)"
$ws.Range("C14").Formula = "=_xlfn.STDEV.S(C2:C11)"
$ws.Range("D14").Formula = "=_xlfn.STDEV.S(D2:D11)"

# --- Page margins (Left/Right 0.75in, Top/Bottom 1in, Header/Footer 0.5in) ---
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# --- Selection moved to D14 ---
$ws.Range("D14").Select()
